$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update card texts (column B) - unchanged rows are left untouched.
# Rows that gain a line-break get wrap text + explicit row height as well.

$ws.Range("B8").Value = "My silence isn't approval.`nI was just not listening."
$ws.Range("B8").WrapText = $true
$ws.Range("B8").EntireRow.RowHeight = 34

$ws.Range("B9").Value = "My video is not frozen.`nI just try not to move."
$ws.Range("B9").WrapText = $true
$ws.Range("B9").EntireRow.RowHeight = 34

$ws.Range("B10").Value = "I'm just trying to stick my mouse pointer in your ear."

$ws.Range("B11").Value = "Your agenda is useless."

$ws.Range("B14").Value = "You already lost me at`n""who sent the invitation for this?"""
$ws.Range("B14").WrapText = $true
$ws.Range("B14").EntireRow.RowHeight = 51

$ws.Range("B15").Value = "I placed a funny photo over your switched-off video."

$ws.Range("B16").Value = "Look there!`nA squirrel!"
$ws.Range("B16").WrapText = $true
$ws.Range("B16").EntireRow.RowHeight = 34

$ws.Range("B17").Value = "Here we are now,`nentertain us!"
$ws.Range("B17").WrapText = $true
$ws.Range("B17").EntireRow.RowHeight = 34

$ws.Range("B18").Value = "All I hear is`n""mi mi mi""."
$ws.Range("B18").WrapText = $true
$ws.Range("B18").EntireRow.RowHeight = 34

# Move the active selection, matching the saved view state.
$ws.Range("B12").Select()
